# Applies the diff: adds a new "Player Info" sheet as the first sheet,
# and renames/updates the MATCH_CARD_LINK column to MATCH_CODE on the
# "ODI Batting" and "ODI Bowling" sheets (storing only the numeric match
# code instead of the full scorecard URL).

$wb = $excel.ActiveWorkbook

# --- 1. Insert a new worksheet named "Player Info" before the first sheet ---
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# Headers
$headers = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $playerInfo.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data row (ID is a numeric-looking value that should stay text, like the
# source "inlineStr" cells, so force a text number format first)
$playerInfo.Cells.Item(2, 1).NumberFormat = "@"
$playerInfo.Cells.Item(2, 1).Value = "5550"
$playerInfo.Cells.Item(2, 2).Value = "George Fredrik Linde"
$playerInfo.Cells.Item(2, 3).Value = "Left Handed"
$playerInfo.Cells.Item(2, 4).Value = "Left Arm Orthodox"

# --- 2. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("D1").Value = "MATCH_CODE"
$odiBatting.Range("D2").NumberFormat = "@"
$odiBatting.Range("D2").Value = "4488"
$odiBatting.Range("D3").NumberFormat = "@"
$odiBatting.Range("D3").Value = "4491"

# --- 3. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$odiBowling = $wb.Worksheets.Item("ODI Bowling")
$odiBowling.Range("B1").Value = "MATCH_CODE"
$odiBowling.Range("B2").NumberFormat = "@"
$odiBowling.Range("B2").Value = "4488"
$odiBowling.Range("B3").NumberFormat = "@"
$odiBowling.Range("B3").Value = "4491"
